# Auto-generated script applying scheduled-runner price/profit updates
# to the Pandaemonium Profits workbook (one worksheet per crafting class).

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value2 = 3200
$ws.Range("I21").Value2 = 3200
$ws.Range("K21").Value2 = 3200
$ws.Range("M21").Value2 = -2732
$ws.Range("H23").Value2 = 3200
$ws.Range("I23").Value2 = 3200
$ws.Range("K23").Value2 = 3200
$ws.Range("M23").Value2 = -2966
$ws.Range("H51").Value2 = 5667
$ws.Range("I51").Value2 = 10001
$ws.Range("J51").Value2 = 3500
$ws.Range("K51").Value2 = 10001
$ws.Range("L51").Value2 = 3500
$ws.Range("M51").Value2 = -9517
$ws.Range("N51").Value2 = -4468
$ws.Range("H64").Value2 = 3661.111
$ws.Range("I64").Value2 = 3541.6667
$ws.Range("J64").Value2 = 3900
$ws.Range("K64").Value2 = 3541.6667
$ws.Range("L64").Value2 = 3900
$ws.Range("M64").Value2 = -3293.6667
$ws.Range("N64").Value2 = -4396
$ws.Range("H67").Value2 = 3661.111
$ws.Range("I67").Value2 = 3541.6667
$ws.Range("J67").Value2 = 3900
$ws.Range("K67").Value2 = 3541.6667
$ws.Range("L67").Value2 = 3900
$ws.Range("M67").Value2 = -2683.6667
$ws.Range("N67").Value2 = -5616
$ws.Range("H74").Value2 = 4312.5
$ws.Range("I74").Value2 = 3750
$ws.Range("J74").Value2 = 6000
$ws.Range("K74").Value2 = 3750
$ws.Range("L74").Value2 = 6000
$ws.Range("M74").Value2 = -2814
$ws.Range("N74").Value2 = -7872
$ws.Range("H75").Value2 = 28000
$ws.Range("J75").Value2 = 0
$ws.Range("L75").Value2 = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value2 = 4312.5
$ws.Range("I77").Value2 = 3750
$ws.Range("J77").Value2 = 6000
$ws.Range("K77").Value2 = 18750
$ws.Range("L77").Value2 = 30000
$ws.Range("M77").Value2 = -14070
$ws.Range("N77").Value2 = -39360
$ws.Range("H78").Value2 = 28000
$ws.Range("J78").Value2 = 0
$ws.Range("L78").Value2 = 0
$ws.Range("N78").ClearContents()
$ws.Range("H135").Value2 = 39475236
$ws.Range("I135").Value2 = 15626352
$ws.Range("J135").Value2 = 166669260
$ws.Range("K135").Value2 = 140637168
$ws.Range("L135").Value2 = 1500023340
$ws.Range("M135").Value2 = -140634633
$ws.Range("N135").Value2 = -1500028410
$ws.Range("H136").Value2 = 64946.668
$ws.Range("J136").Value2 = 64946.668
$ws.Range("L136").Value2 = 64946.668
$ws.Range("N136").Value2 = -75146.66800000001
$ws.Range("H137").Value2 = 2666.2
$ws.Range("I137").Value2 = 1748.0286
$ws.Range("J137").Value2 = 4273
$ws.Range("K137").Value2 = 5244.085800000001
$ws.Range("L137").Value2 = 12819
$ws.Range("M137").Value2 = -2694.085800000001
$ws.Range("N137").Value2 = -17919
$ws.Range("H138").Value2 = 3585.4724
$ws.Range("I138").Value2 = 1681.2941
$ws.Range("J138").Value2 = 4721.2983
$ws.Range("K138").Value2 = 5043.8823
$ws.Range("L138").Value2 = 14163.8949
$ws.Range("M138").Value2 = 96.11769999999979
$ws.Range("N138").Value2 = -24443.8949
$ws.Range("H140").Value2 = 102440.766
$ws.Range("J140").Value2 = 102440.766
$ws.Range("L140").Value2 = 102440.766
$ws.Range("N140").Value2 = -112800.766
$ws.Range("H141").Value2 = 3173.3333
$ws.Range("I141").Value2 = 2710.6667
$ws.Range("J141").Value2 = 3944.4443
$ws.Range("K141").Value2 = 8132.000100000001
$ws.Range("L141").Value2 = 11833.3329
$ws.Range("M141").Value2 = -2952.000100000001
$ws.Range("N141").Value2 = -22193.3329

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value2 = 35402.168
$ws.Range("J23").Value2 = 26481.4
$ws.Range("L23").Value2 = 26481.4
$ws.Range("N23").Value2 = -26999.4
$ws.Range("H32").Value2 = 9913.192999999999
$ws.Range("I32").Value2 = 9174.076999999999
$ws.Range("J32").Value2 = 17600
$ws.Range("K32").Value2 = 9174.076999999999
$ws.Range("L32").Value2 = 17600
$ws.Range("M32").Value2 = -8887.076999999999
$ws.Range("N32").Value2 = -18174
$ws.Range("H33").Value2 = 12966.667
$ws.Range("I33").Value2 = 6000
$ws.Range("J33").Value2 = 19933.334
$ws.Range("K33").Value2 = 6000
$ws.Range("L33").Value2 = 19933.334
$ws.Range("M33").Value2 = -5671
$ws.Range("N33").Value2 = -20591.334
$ws.Range("H55").Value2 = 125051624
$ws.Range("J55").Value2 = 125051624
$ws.Range("L55").Value2 = 125051624
$ws.Range("N55").Value2 = -125052254
$ws.Range("H88").Value2 = 5113.8
$ws.Range("I88").Value2 = 8459
$ws.Range("J88").Value2 = 2186.75
$ws.Range("K88").Value2 = 8459
$ws.Range("L88").Value2 = 2186.75
$ws.Range("M88").Value2 = -8053
$ws.Range("N88").Value2 = -2998.75
$ws.Range("H91").Value2 = 5113.8
$ws.Range("I91").Value2 = 8459
$ws.Range("J91").Value2 = 2186.75
$ws.Range("K91").Value2 = 8459
$ws.Range("L91").Value2 = 2186.75
$ws.Range("M91").Value2 = -7055
$ws.Range("N91").Value2 = -4994.75
$ws.Range("H110").Value2 = 1965.1904
$ws.Range("I110").Value2 = 1816.8125
$ws.Range("J110").Value2 = 2440
$ws.Range("K110").Value2 = 1816.8125
$ws.Range("L110").Value2 = 2440
$ws.Range("M110").Value2 = 228.1875
$ws.Range("N110").Value2 = -6530
$ws.Range("H132").Value2 = 2360.8235
$ws.Range("I132").Value2 = 1318.625
$ws.Range("K132").Value2 = 3955.875
$ws.Range("M132").Value2 = -1425.875

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value2 = 37000
$ws.Range("J35").Value2 = 37000
$ws.Range("L35").Value2 = 37000
$ws.Range("N35").Value2 = -37620
$ws.Range("H94").Value2 = 2336
$ws.Range("I94").Value2 = 1903.3334
$ws.Range("K94").Value2 = 1903.3334
$ws.Range("M94").Value2 = -1452.3334
$ws.Range("H105").Value2 = 4567.381
$ws.Range("I105").Value2 = 4448.727
$ws.Range("J105").Value2 = 5002.4443
$ws.Range("K105").Value2 = 4448.727
$ws.Range("L105").Value2 = 5002.4443
$ws.Range("M105").Value2 = -2701.727
$ws.Range("N105").Value2 = -8496.444299999999
$ws.Range("H107").Value2 = 2500
$ws.Range("I107").Value2 = 2000
$ws.Range("J107").Value2 = 4000
$ws.Range("K107").Value2 = 2000
$ws.Range("L107").Value2 = 4000
$ws.Range("M107").Value2 = -80
$ws.Range("N107").Value2 = -7840

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3732.651
$ws.Range("I31").Value2 = 2046.7084
$ws.Range("J31").Value2 = 9127.666999999999
$ws.Range("K31").Value2 = 2046.7084
$ws.Range("L31").Value2 = 9127.666999999999
$ws.Range("M31").Value2 = -1751.7084
$ws.Range("N31").Value2 = -9717.666999999999
$ws.Range("H34").Value2 = 3732.651
$ws.Range("I34").Value2 = 2046.7084
$ws.Range("J34").Value2 = 9127.666999999999
$ws.Range("K34").Value2 = 2046.7084
$ws.Range("L34").Value2 = 9127.666999999999
$ws.Range("M34").Value2 = -1844.7084
$ws.Range("N34").Value2 = -9531.666999999999
$ws.Range("H41").Value2 = 9999.5
$ws.Range("I41").Value2 = 9999.5
$ws.Range("K41").Value2 = 9999.5
$ws.Range("M41").Value2 = -9571.5
$ws.Range("H60").Value2 = 17000
$ws.Range("J60").Value2 = 31000
$ws.Range("L60").Value2 = 31000
$ws.Range("N60").Value2 = -32022
$ws.Range("H107").Value2 = 1572.6
$ws.Range("I107").Value2 = 1545.8889
$ws.Range("J107").Value2 = 1813
$ws.Range("K107").Value2 = 1545.8889
$ws.Range("L107").Value2 = 1813
$ws.Range("M107").Value2 = 374.1111000000001
$ws.Range("N107").Value2 = -5653

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 11114356
$ws.Range("I5").Value2 = 626.6070999999999
$ws.Range("J5").Value2 = 29419322
$ws.Range("K5").Value2 = 1879.8213
$ws.Range("L5").Value2 = 88257966
$ws.Range("M5").Value2 = -1767.8213
$ws.Range("N5").Value2 = -88258190
$ws.Range("H75").Value2 = 1349.5714
$ws.Range("J75").Value2 = 1349.5714
$ws.Range("L75").Value2 = 4048.7142
$ws.Range("N75").Value2 = -6044.7142
$ws.Range("H78").Value2 = 1349.5714
$ws.Range("J78").Value2 = 1349.5714
$ws.Range("L78").Value2 = 12146.1426
$ws.Range("N78").Value2 = -22130.1426
$ws.Range("H106").Value2 = 5219.846
$ws.Range("J106").Value2 = 5219.846
$ws.Range("L106").Value2 = 15659.538
$ws.Range("N106").Value2 = -17551.538
$ws.Range("H131").Value2 = 11501.632
$ws.Range("I131").Value2 = 430.0635
$ws.Range("J131").Value2 = 40564.5
$ws.Range("K131").Value2 = 1290.1905
$ws.Range("L131").Value2 = 121693.5
$ws.Range("M131").Value2 = 3749.8095
$ws.Range("N131").Value2 = -131773.5
$ws.Range("H135").Value2 = 11114356
$ws.Range("I135").Value2 = 626.6070999999999
$ws.Range("J135").Value2 = 29419322
$ws.Range("K135").Value2 = 5639.4639
$ws.Range("L135").Value2 = 264773898
$ws.Range("M135").Value2 = -3104.4639
$ws.Range("N135").Value2 = -264778968
$ws.Range("H137").Value2 = 41670830
$ws.Range("J137").Value2 = 6196.2
$ws.Range("L137").Value2 = 18588.6
$ws.Range("N137").Value2 = -28788.6

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 23104.604
$ws.Range("I132").Value2 = 46760.137
$ws.Range("J132").Value2 = 3088.3845
$ws.Range("K132").Value2 = 140280.411
$ws.Range("L132").Value2 = 9265.1535
$ws.Range("M132").Value2 = -137750.411
$ws.Range("N132").Value2 = -14325.1535
$ws.Range("H135").Value2 = 39960
$ws.Range("J135").Value2 = 39960
$ws.Range("L135").Value2 = 39960
$ws.Range("N135").Value2 = -50100

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 3001
$ws.Range("I100").Value2 = 1001
$ws.Range("J100").Value2 = 4334.3335
$ws.Range("K100").Value2 = 2002
$ws.Range("L100").Value2 = 8668.666999999999
$ws.Range("M100").Value2 = -1461
$ws.Range("N100").Value2 = -9750.666999999999
$ws.Range("H107").Value2 = 1161.16
$ws.Range("I107").Value2 = 696.2778
$ws.Range("J107").Value2 = 2356.5715
$ws.Range("K107").Value2 = 2088.8334
$ws.Range("L107").Value2 = 7069.7145
$ws.Range("M107").Value2 = -168.8334
$ws.Range("N107").Value2 = -10909.7145
$ws.Range("H132").Value2 = 2540.4666
$ws.Range("I132").Value2 = 2100.861
$ws.Range("J132").Value2 = 4298.8887
$ws.Range("K132").Value2 = 6302.583
$ws.Range("L132").Value2 = 12896.6661
$ws.Range("M132").Value2 = -3772.583
$ws.Range("N132").Value2 = -17956.6661
